# Update "want to go" counts (column F) and one sold-out price cell (G9)
# across the workbook's sheets, per the gh-pages regeneration commit.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value  = 1506
$ws1.Range("F9").Value  = 38851
$ws1.Range("G9").Value  = 0
$ws1.Range("F10").Value = 8165
$ws1.Range("F18").Value = 625
$ws1.Range("F19").Value = 41
$ws1.Range("F21").Value = 540
$ws1.Range("F23").Value = 1012
$ws1.Range("F31").Value = 366
$ws1.Range("F33").Value = 807
$ws1.Range("F34").Value = 351
$ws1.Range("F36").Value = 219
$ws1.Range("F37").Value = 955

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value  = 6
$ws2.Range("F9").Value  = 5
$ws2.Range("F15").Value = 0

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 1506
$ws4.Range("F16").Value = 8165
$ws4.Range("F28").Value = 540
$ws4.Range("F29").Value = 1012
$ws4.Range("F37").Value = 366
$ws4.Range("F39").Value = 807
